$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 94: DfT Group / Rail / Rail / HSMRPG / RDP / Rail
$ws.Range("B94").Value = "Rail"
$ws.Range("C94").Value = "Rail"
$ws.Range("D94").Value = "HSMRPG"
$ws.Range("E94").Value = "RDP"
$ws.Range("F94").Value = "Rail"

# Row 95: IPDC approval point / FBC / FBC / SOBC / OBC / OBC
$ws.Range("B95").Value = "FBC"
$ws.Range("C95").Value = "FBC"
$ws.Range("D95").Value = "SOBC"
$ws.Range("E95").Value = "OBC"
$ws.Range("F95").Value = "OBC"

# Update selection to reflect the last active cell after edits
$ws.Range("D107").Select()
